$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 374594.78
$ws.Range("J17").Value = 387463.8
$ws.Range("L17").Value = 1162391.4
$ws.Range("N17").Value = -1162727.4
# Row 51
$ws.Range("H51").Value = 4373.339
$ws.Range("I51").Value = 3269.5557
$ws.Range("K51").Value = 3269.5557
$ws.Range("M51").Value = -2785.5557
# Row 103
$ws.Range("H103").Value = 520.8
$ws.Range("J103").Value = 1000
$ws.Range("L103").Value = 3000
$ws.Range("N103").Value = -4172
# Row 112
$ws.Range("H112").Value = 41886.8
$ws.Range("I112").Value = 336000
$ws.Range("J112").Value = 1780.4546
$ws.Range("K112").Value = 1008000
$ws.Range("L112").Value = 5341.3638
$ws.Range("M112").Value = -1006892
$ws.Range("N112").Value = -7557.3638
# Row 125
$ws.Range("H125").Value = 1454.8
$ws.Range("I125").Value = 1167.3334
$ws.Range("K125").Value = 10506.0006
$ws.Range("M125").Value = -8046.000599999999
# Row 138
$ws.Range("H138").Value = 1695.9744
$ws.Range("I138").Value = 1313.069
$ws.Range("K138").Value = 3939.207
$ws.Range("M138").Value = 1200.793

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 2998.7334
$ws.Range("I5").Value = 280.72726
$ws.Range("K5").Value = 280.72726
$ws.Range("M5").Value = -168.72726
# Row 25
$ws.Range("H25").Value = 2723
$ws.Range("I25").Value = 2723
$ws.Range("K25").Value = 2723
$ws.Range("M25").Value = -2321
# Row 61
$ws.Range("H61").Value = 1291.5883
$ws.Range("I61").Value = 1128.5807
$ws.Range("K61").Value = 1128.5807
$ws.Range("M61").Value = -916.5807
# Row 74
$ws.Range("H74").Value = 2656.125
$ws.Range("I74").Value = 2589.9656
$ws.Range("J74").Value = 3295.6667
$ws.Range("K74").Value = 2589.9656
$ws.Range("L74").Value = 3295.6667
$ws.Range("M74").Value = -1715.9656
$ws.Range("N74").Value = -5043.6667
# Row 77
$ws.Range("H77").Value = 2656.125
$ws.Range("I77").Value = 2589.9656
$ws.Range("J77").Value = 3295.6667
$ws.Range("K77").Value = 12949.828
$ws.Range("L77").Value = 16478.3335
$ws.Range("M77").Value = -8581.828
$ws.Range("N77").Value = -25214.3335
# Row 132
$ws.Range("H132").Value = 938.9167
$ws.Range("I132").Value = 872.35486
$ws.Range("K132").Value = 2617.06458
$ws.Range("M132").Value = -87.06458000000021
# Row 136
$ws.Range("H136").Value = 1291.5883
$ws.Range("I136").Value = 1128.5807
$ws.Range("K136").Value = 3385.7421
$ws.Range("M136").Value = -835.7420999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 2998.7334
$ws.Range("I4").Value = 280.72726
$ws.Range("K4").Value = 280.72726
$ws.Range("M4").Value = -165.72726
# Row 15
$ws.Range("H15").Value = 607
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 607
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 607
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -1061
# Row 68
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
# Row 71
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
# Row 107
$ws.Range("H107").Value = 41667964
$ws.Range("I107").Value = 1061.2222
$ws.Range("K107").Value = 1061.2222
$ws.Range("M107").Value = 858.7778000000001
# Row 134
$ws.Range("H134").Value = 799.6667
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
# Row 137
$ws.Range("H137").Value = 55000
$ws.Range("J137").Value = 55000
$ws.Range("L137").Value = 55000
$ws.Range("N137").Value = -65200
# Row 138
$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9429.117
$ws.Range("I31").Value = 21964
$ws.Range("J31").Value = 2591.9092
$ws.Range("K31").Value = 21964
$ws.Range("L31").Value = 2591.9092
$ws.Range("M31").Value = -21669
$ws.Range("N31").Value = -3181.9092
# Row 34
$ws.Range("H34").Value = 9429.117
$ws.Range("I34").Value = 21964
$ws.Range("J34").Value = 2591.9092
$ws.Range("K34").Value = 21964
$ws.Range("L34").Value = 2591.9092
$ws.Range("M34").Value = -21762
$ws.Range("N34").Value = -2995.9092
# Row 132
$ws.Range("H132").Value = 3186.6
$ws.Range("I132").Value = 2915.7585
$ws.Range("J132").Value = 4495.6665
$ws.Range("K132").Value = 8747.2755
$ws.Range("L132").Value = 13486.9995
$ws.Range("M132").Value = -6217.2755
$ws.Range("N132").Value = -18546.9995
# Row 134
$ws.Range("H134").Value = 1670.8235
$ws.Range("I134").Value = 1427.0667
$ws.Range("J134").Value = 3499
$ws.Range("K134").Value = 4281.2001
$ws.Range("L134").Value = 10497
$ws.Range("M134").Value = -1746.2001
$ws.Range("N134").Value = -15567

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 16620847
$ws.Range("I4").Value = 36849936
$ws.Range("J4").Value = 4096.7856
$ws.Range("K4").Value = 110549808
$ws.Range("L4").Value = 12290.3568
$ws.Range("M4").Value = -110549696
$ws.Range("N4").Value = -12514.3568
# Row 20
$ws.Range("H20").Value = 500
$ws.Range("I20").Value = 500
$ws.Range("K20").Value = 1500
$ws.Range("M20").Value = -1273
# Row 36
$ws.Range("H36").Value = 999
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
# Row 55
$ws.Range("H55").Value = 1522
$ws.Range("I55").Value = 873
$ws.Range("J55").Value = 2333.25
$ws.Range("K55").Value = 2619
$ws.Range("L55").Value = 6999.75
$ws.Range("M55").Value = -2442
$ws.Range("N55").Value = -7353.75
# Row 58
$ws.Range("H58").Value = 1833
$ws.Range("I58").Value = 1499
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 4497
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -4369
$ws.Range("N58").Value = -6256
# Row 81
$ws.Range("H81").Value = 8613.5
$ws.Range("J81").Value = 8613.5
$ws.Range("L81").Value = 25840.5
$ws.Range("N81").Value = -28086.5
# Row 84
$ws.Range("H84").Value = 8613.5
$ws.Range("J84").Value = 8613.5
$ws.Range("L84").Value = 77521.5
$ws.Range("N84").Value = -88753.5
# Row 109
$ws.Range("H109").Value = 365.8889
$ws.Range("I109").Value = 365.8889
$ws.Range("K109").Value = 1097.6667
$ws.Range("M109").Value = -57.66669999999999
# Row 122
$ws.Range("H122").Value = 963
$ws.Range("I122").Value = 969.5
$ws.Range("J122").Value = 950
$ws.Range("K122").Value = 8725.5
$ws.Range("L122").Value = 8550
$ws.Range("M122").Value = -6275.5
$ws.Range("N122").Value = -13450
# Row 134
$ws.Range("H134").Value = 2009.7333
$ws.Range("I134").Value = 1439
$ws.Range("K134").Value = 4317
$ws.Range("M134").Value = 753
# Row 140
$ws.Range("H140").Value = 2998.3635
$ws.Range("I140").Value = 1054.75
$ws.Range("K140").Value = 3164.25
$ws.Range("M140").Value = 2015.75

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2756.1404
$ws.Range("I132").Value = 2151.25
$ws.Range("K132").Value = 6453.75
$ws.Range("M132").Value = -3923.75

$ws = $wb.Worksheets.Item("LTW")
# Row 6
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
# Row 43
$ws.Range("H43").Value = 6668750
$ws.Range("I43").Value = 4170833.2
$ws.Range("J43").Value = 9166667
$ws.Range("K43").Value = 4170833.2
$ws.Range("L43").Value = 9166667
$ws.Range("M43").Value = -4170640.2
$ws.Range("N43").Value = -9167053
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1514.557
$ws.Range("I132").Value = 1344.9546
$ws.Range("K132").Value = 4034.8638
$ws.Range("M132").Value = -1504.8638
# Row 136
$ws.Range("H136").Value = 1964.9592
$ws.Range("I136").Value = 967.7368
$ws.Range("K136").Value = 2903.2104
$ws.Range("M136").Value = -353.2103999999999
